$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$ws.Range("B43").Value = "p<=15" + $nl + "p is definied in executive decree 37121"
$ws.Range("B44").Value = "15>p=<30" + $nl + "p is definied in executive decree 37121"
$ws.Range("B45").Value = "30<p<=100" + $nl + "p is definied in executive decree 37121"
$ws.Range("B46").Value = "p > 100" + $nl + "p is definied in executive decree 37121"
